# The commit renames the documented helper function "str_split_n()" to
# "str_split_i()" in the stringr cheat-sheet style table on slide 1
# (inside the text box that documents str_c / str_flatten / str_dup /
# str_split_fixed / str_glue / str_glue_data), updating both the inline
# code reference ("str_split_n()") and the following prose ("nth
# substring" -> "ith substring").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape by its distinctive name/content instead of a hard-coded
# index, so the script is resilient to shape re-ordering.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        $t = $candidate.TextFrame.TextRange.Text
        if ($t.IndexOf("str_split_n") -ge 0) {
            $shp = $candidate
            break
        }
    }
}

if ($shp -ne $null) {
    $tr = $shp.TextFrame.TextRange

    # 1) "str_split_n(" -> "str_split_i("
    $full = $tr.Text
    $idx = $full.IndexOf("str_split_n")
    if ($idx -ge 0) {
        $rng = $tr.Characters($idx + 1, "str_split_n".Length)
        $rng.Text = "str_split_i"
    }

    # 2) "nth substring" -> "ith substring"
    $full = $tr.Text
    $idx = $full.IndexOf("nth substring")
    if ($idx -ge 0) {
        $rng = $tr.Characters($idx + 1, "nth".Length)
        $rng.Text = "ith"
    }
}
